$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-7
$data = @(
    @("Charlotte Hornets vs Phoenix Suns", "15-03-2024", "Charlotte"),
    @("Detroit Pistons vs Miami Heat", "15-03-2024", "Detroit"),
    @("Toronto Raptors vs Orlando Magic", "15-03-2024", "Toronto"),
    @("New Orleans Pelicans vs LA Clippers", "15-03-2024", "New Orleans"),
    @("San Antonio Spurs vs Denver Nuggets", "15-03-2024", "Austin"),
    @("Utah Jazz vs Atlanta Hawks", "15-03-2024", "Salt Lake City")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# Remove rows 8-11 (now obsolete)
$ws.Range("A8:C11").ClearContents()
